$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "Integer min" value for rule R30 (row 10) from 18 to 1
$ws.Range("C10").Value = 1
